# EnemyImbuePresets_MenuMock.xlsx - "Menu Layout" sheet edits
# Refine enemy type presets and simplify logging controls

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Layout")

# --- Row 94: Diagnostics/Log Level -> Advanced/Basic Logs (dropdown -> toggle) ---
$ws.Cells.Item(94, 2).Value = "Advanced"
$ws.Cells.Item(94, 3).Value = "Basic Logs"
$ws.Cells.Item(94, 4).Value = "Toggle"
$ws.Cells.Item(94, 5).Value = "On/Off"
$ws.Cells.Item(94, 6).Value = "On"
$ws.Cells.Item(94, 8).Value = "General informational logging"

# --- Row 95: Diagnostics/Imbue Update Interval -> Advanced/Diagnostics Logs (toggle) ---
$ws.Cells.Item(95, 2).Value = "Advanced"
$ws.Cells.Item(95, 3).Value = "Diagnostics Logs"
$ws.Cells.Item(95, 4).Value = "Toggle"
$ws.Cells.Item(95, 5).Value = "On/Off"
$ws.Cells.Item(95, 6).Value = "Off"
$ws.Cells.Item(95, 8).Value = "Deeper troubleshooting logs"

# --- Row 96: Diagnostics/Enemy Rescan Interval -> Advanced/Verbose Logs (toggle) ---
$ws.Cells.Item(96, 2).Value = "Advanced"
$ws.Cells.Item(96, 3).Value = "Verbose Logs"
$ws.Cells.Item(96, 4).Value = "Toggle"
$ws.Cells.Item(96, 5).Value = "On/Off"
$ws.Cells.Item(96, 6).Value = "Off"
$ws.Cells.Item(96, 8).Value = "High-volume per-creature logs"

# --- Row 97: Diagnostics/Dump Factions -> Advanced/Session Diagnostics (toggle) ---
$ws.Cells.Item(97, 2).Value = "Advanced"
$ws.Cells.Item(97, 3).Value = "Session Diagnostics"
$ws.Cells.Item(97, 4).Value = "Toggle"
$ws.Cells.Item(97, 5).Value = "On/Off"
$ws.Cells.Item(97, 6).Value = "Off"
$ws.Cells.Item(97, 8).Value = "Structured session summary logs"

# --- Row 98: Diagnostics/Dump Wave-Faction Map -> Advanced/Imbue Update Interval (dropdown) ---
$ws.Cells.Item(98, 2).Value = "Advanced"
$ws.Cells.Item(98, 3).Value = "Imbue Update Interval"
$ws.Cells.Item(98, 4).Value = "Dropdown"
$ws.Cells.Item(98, 5).Value = "0.05s..1.00s"
$ws.Cells.Item(98, 6).Value = "0.25s"
$ws.Cells.Item(98, 8).Value = "Performance/response tradeoff"

# --- Row 99: Diagnostics/Dump State -> Advanced/Enemy Rescan Interval (dropdown) ---
$ws.Cells.Item(99, 2).Value = "Advanced"
$ws.Cells.Item(99, 3).Value = "Enemy Rescan Interval"
$ws.Cells.Item(99, 4).Value = "Dropdown"
$ws.Cells.Item(99, 5).Value = "0.50s..5.00s"
$ws.Cells.Item(99, 6).Value = "2.00s"
$ws.Cells.Item(99, 8).Value = "Tracking refresh interval"

# --- Row 100: Diagnostics/Dump Enemy Type Detection -> Advanced/Force Reapply (button, unchanged D/E/F) ---
$ws.Cells.Item(100, 2).Value = "Advanced"
$ws.Cells.Item(100, 3).Value = "Force Reapply"
$ws.Cells.Item(100, 8).Value = "One-shot reapply action"

# --- Row 101 (old Force Reapply) is removed entirely; rows shift up ---
$ws.Rows(101).Delete()

# --- Fix up sheet dimension-dependent refs: autofilter now A1:H100 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:H100").AutoFilter()

# --- Update the workbook-level _FilterDatabase defined name for this sheet ---
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = '=''Menu Layout''!$A$1:$H$100'
